# repull data, push all data, mean calculation
# Update the "dSF" (F) column values to reflect the repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 0
    4  = -1
    6  = 0
    9  = -5
    10 = -7
    11 = 3
    12 = -5
    13 = -7
    15 = -5
    16 = -6
    18 = -2
    19 = -3
    20 = -7
    22 = 0
    23 = -6
    25 = 4
    26 = -5
    27 = 5
    28 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
